$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 and J1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting (bold, centered, bordered style) from H1 onto I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new I/J column data for rows 2-14
$data = @(
    @(6, 9),
    @(3, 6),
    @(5, 7),
    @(9, 9),
    @(8, 8),
    @(6, 9),
    @(1, 5),
    @(1, 4),
    @(1, 4),
    @(8, 9),
    @(8, 9),
    @(7, 7),
    @(1, 2)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
